$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 21 (shifts existing row 21 and below down by one)
$ws.Rows.Item(21).Insert()

# Fill the newly inserted row with the new city data.
$ws.Range("A21").Value = "Kővágóörs"

# Coordinates are stored as TEXT in this sheet (not numbers). A direct
# ws.Range(...).Value = "46.84921" assignment would be auto-converted to
# a number by Excel. Instead, build the text via a throwaway formula
# cell (whose result is a genuine text value) and paste-special just the
# values into the target cells - this keeps them as text without
# introducing any new cell styles/number formats.
$scratch = $ws.Range("Z1")
$scratch.Formula = '="46.84921"'
$scratch.Copy()
$ws.Range("B21").PasteSpecial(-4163)  # xlPasteValues

$scratch.Formula = '="17.59911"'
$scratch.Copy()
$ws.Range("C21").PasteSpecial(-4163)  # xlPasteValues

$scratch.ClearContents()

# Update the selection to match the saved workbook view state
$ws.Range("A21").Select()
